# Updated cryptos list on Thu May 11 09:16:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.447.75"
$ws.Range("E2").Value = "  -0.61%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.823.97"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'312.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.01%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4230"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.55%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.37%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.8593"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'20.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.802.57"
$ws.Range("E12").Value = "  -1.31%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.391"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'6.469"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.06930"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  +0.00%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'80.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000008889"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.08%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'15.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.724.48"
$ws.Range("E21").Value = "  +0.60%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.138"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.55%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  +5.50%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.001.33"
$ws.Range("E24").Value = "  -1.90%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'1.985"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'154.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'18.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'5.155"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "'114.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.46%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "'1.793"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.66%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08821"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32 and 33 swap: HuobiToken <-> ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7486"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.966"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'4.528"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "

# Row 36 - Frax
$ws.Range("E36").Value = "  +0.05%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  -1.68%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "'0.05280"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.01918"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.770"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.5055"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "

# Row 42 - Algorand
$ws.Range("D42").Value = "'0.1639"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'6.452"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'8.314"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "'10.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.62%  "

# Row 46 - Quant
$ws.Range("D46").Value = "'105.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "

# Row 47 - Cronos
$ws.Range("D47").Value = "'0.06444"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.43%  "

# Row 48 - Decentraland
$ws.Range("D48").Value = "'0.4670"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.07%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  -0.02%  "

# Row 50 - NEARProtocol
$ws.Range("E50").Value = "  -1.06%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'63.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
